$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 556.3421
$ws.Range("J17").Value = 556.3421
$ws.Range("L17").Value = 1669.0263
$ws.Range("N17").Value = -2005.0263

$ws.Range("H33").Value = 854.44446
$ws.Range("I33").Value = 711.375
$ws.Range("J33").Value = 1999
$ws.Range("K33").Value = 711.375
$ws.Range("L33").Value = 1999
$ws.Range("M33").Value = -482.375
$ws.Range("N33").Value = -2457

$ws.Range("H40").Value = 6264.6665
$ws.Range("I40").Value = 5358
$ws.Range("J40").Value = 6990
$ws.Range("K40").Value = 5358
$ws.Range("L40").Value = 6990
$ws.Range("M40").Value = -5183
$ws.Range("N40").Value = -7340

$ws.Range("H86").Value = 3301.2
$ws.Range("I86").Value = 3837.6667
$ws.Range("K86").Value = 3837.6667
$ws.Range("M86").Value = -2714.6667

$ws.Range("H89").Value = 3301.2
$ws.Range("I89").Value = 3837.6667
$ws.Range("K89").Value = 19188.3335
$ws.Range("M89").Value = -13572.3335

$ws.Range("H92").Value = 2271.625
$ws.Range("I92").Value = 2094.8333
$ws.Range("J92").Value = 2802
$ws.Range("K92").Value = 2094.8333
$ws.Range("L92").Value = 2802
$ws.Range("M92").Value = -846.8332999999998
$ws.Range("N92").Value = -5298

$ws.Range("H107").Value = 2596.5
$ws.Range("I107").Value = 1403
$ws.Range("J107").Value = 5381.3335
$ws.Range("K107").Value = 1403
$ws.Range("L107").Value = 5381.3335
$ws.Range("M107").Value = 517
$ws.Range("N107").Value = -9221.333500000001

$ws.Range("H112").Value = 3191.353
$ws.Range("J112").Value = 2016.8667
$ws.Range("L112").Value = 6050.6001
$ws.Range("N112").Value = -8266.6001

$ws.Range("H125").Value = 10002
$ws.Range("I125").Value = 10002
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 90018
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -87558
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 498543.06
$ws.Range("I132").Value = 615422.4399999999
$ws.Range("K132").Value = 1846267.32
$ws.Range("M132").Value = -1843737.32

$ws.Range("H136").Value = 144854
$ws.Range("J136").Value = 148999
$ws.Range("L136").Value = 148999
$ws.Range("N136").Value = -159199

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5060.615
$ws.Range("I74").Value = 2699.8
$ws.Range("K74").Value = 2699.8
$ws.Range("M74").Value = -1825.8

$ws.Range("H77").Value = 5060.615
$ws.Range("I77").Value = 2699.8
$ws.Range("K77").Value = 13499
$ws.Range("M77").Value = -9131

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 1999
$ws.Range("J19").Value = 2499
$ws.Range("L19").Value = 2499
$ws.Range("N19").Value = -2845

$ws.Range("H20").Value = 8405583
$ws.Range("I20").Value = 20409458
$ws.Range("J20").Value = 2870.5
$ws.Range("K20").Value = 20409458
$ws.Range("L20").Value = 2870.5
$ws.Range("M20").Value = -20409211
$ws.Range("N20").Value = -3364.5

$ws.Range("H22").Value = 2194.1428
$ws.Range("I22").Value = 309.1111
$ws.Range("K22").Value = 309.1111
$ws.Range("M22").Value = -136.1111

$ws.Range("H99").Value = 6328.28
$ws.Range("I99").Value = 2877.4614
$ws.Range("J99").Value = 10066.667
$ws.Range("K99").Value = 2877.4614
$ws.Range("L99").Value = 10066.667
$ws.Range("M99").Value = -1379.4614
$ws.Range("N99").Value = -13062.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6750.615
$ws.Range("I132").Value = 5931.5293
$ws.Range("J132").Value = 8297.777
$ws.Range("K132").Value = 17794.5879
$ws.Range("L132").Value = 24893.331
$ws.Range("M132").Value = -15264.5879
$ws.Range("N132").Value = -29953.331

$ws.Range("H134").Value = 31259730
$ws.Range("I134").Value = 125005940
$ws.Range("K134").Value = 375017820
$ws.Range("M134").Value = -375015285

$ws.Range("H140").Value = 158119.33
$ws.Range("J140").Value = 158119.33
$ws.Range("L140").Value = 158119.33
$ws.Range("N140").Value = -168479.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12088.55
$ws.Range("I3").Value = 3905.9167
$ws.Range("K3").Value = 11717.7501
$ws.Range("M3").Value = -11605.7501

$ws.Range("H4").Value = 65597.89999999999
$ws.Range("I4").Value = 442.83334
$ws.Range("K4").Value = 1328.50002
$ws.Range("M4").Value = -1216.50002

$ws.Range("H63").Value = 12717.9375
$ws.Range("I63").Value = 7109.778
$ws.Range("J63").Value = 19928.428
$ws.Range("K63").Value = 21329.334
$ws.Range("L63").Value = 59785.284
$ws.Range("M63").Value = -20580.334
$ws.Range("N63").Value = -61283.284

$ws.Range("H66").Value = 12717.9375
$ws.Range("I66").Value = 7109.778
$ws.Range("J66").Value = 19928.428
$ws.Range("K66").Value = 63988.002
$ws.Range("L66").Value = 179355.852
$ws.Range("M66").Value = -60244.002
$ws.Range("N66").Value = -186843.852

$ws.Range("H69").Value = 12000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 7245.1665
$ws.Range("I70").Value = 5618
$ws.Range("K70").Value = 16854
$ws.Range("M70").Value = -16539

$ws.Range("H72").Value = 12000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 7245.1665
$ws.Range("I73").Value = 5618
$ws.Range("K73").Value = 16854
$ws.Range("M73").Value = -15762

$ws.Range("H76").Value = 24900
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 24900
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H81").Value = 1930.6666
$ws.Range("I81").Value = 1930.6666
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5791.9998
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4668.9998
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1930.6666
$ws.Range("I84").Value = 1930.6666
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 17375.9994
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -11759.9994
$ws.Range("N84").ClearContents()

$ws.Range("H98").Value = 2298.8572
$ws.Range("J98").Value = 1848.6666
$ws.Range("L98").Value = 5545.9998
$ws.Range("N98").Value = -8541.9998

$ws.Range("H117").Value = 1880.3077
$ws.Range("J117").Value = 843
$ws.Range("L117").Value = 2529
$ws.Range("N117").Value = -9413

$ws.Range("H122").Value = 116806.8
$ws.Range("J122").Value = 116806.8
$ws.Range("L122").Value = 1051261.2
$ws.Range("N122").Value = -1056161.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12209
$ws.Range("I122").Value = 7835
$ws.Range("J122").Value = 16583
$ws.Range("K122").Value = 23505
$ws.Range("L122").Value = 49749
$ws.Range("M122").Value = -21055
$ws.Range("N122").Value = -54649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10531.2
$ws.Range("I7").Value = 9260.333000000001
$ws.Range("K7").Value = 9260.333000000001
$ws.Range("M7").Value = -9148.333000000001

$ws.Range("H100").Value = 2023.75
$ws.Range("I100").Value = 1789.5
$ws.Range("K100").Value = 1789.5
$ws.Range("M100").Value = -1248.5

$ws.Range("H122").Value = 5466.5
$ws.Range("I122").Value = 3332.6667
$ws.Range("K122").Value = 9998.000100000001
$ws.Range("M122").Value = -7548.000100000001

$ws.Range("H126").Value = 10531.2
$ws.Range("I126").Value = 9260.333000000001
$ws.Range("K126").Value = 27780.999
$ws.Range("M126").Value = -25310.999

$ws.Range("H132").Value = 2877.5
$ws.Range("I132").Value = 2326.4285
$ws.Range("J132").Value = 4163.3335
$ws.Range("K132").Value = 6979.2855
$ws.Range("L132").Value = 12490.0005
$ws.Range("M132").Value = -4449.2855
$ws.Range("N132").Value = -17550.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3180.7896
$ws.Range("I122").Value = 2887.4443
$ws.Range("J122").Value = 3444.8
$ws.Range("K122").Value = 8662.332900000001
$ws.Range("L122").Value = 10334.4
$ws.Range("M122").Value = -6212.332900000001
$ws.Range("N122").Value = -15234.4

$ws.Range("H132").Value = 7138.1304
$ws.Range("I132").Value = 3525.2727
$ws.Range("K132").Value = 10575.8181
$ws.Range("M132").Value = -8045.8181

Write-Output "applied all changes"